# Applies the "break out stock.yaml completed" update to the "day" sheet:
#   1) D526:D533 (existing rows) were stored as text (inlineStr) BSE codes;
#      re-enter them as plain numbers (same values, numeric type).
#   2) Append 11 new rows (534:544) of freshly scraped stock data, with the
#      bsecode column (D) kept as text, matching the source feed's format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# --- Step 1: re-key D526:D533 as numeric values (was inlineStr text) ---
$bseExisting = @(540699, 500495, 500480, 532286, 500228, 531642, 533278, 500113)
for ($i = 0; $i -lt $bseExisting.Length; $i++) {
    $row = 526 + $i
    $ws.Cells.Item($row, 4).Value = $bseExisting[$i]
}

# --- Step 2: append the new rows 534:544 ---
function Set-StockRow {
    param($ws, $row, $sr, $nsecode, $name, $bsecode, $perChg, $close, $volume, $timeframe, $dateTime)
    $ws.Cells.Item($row, 1).Value = $sr
    $ws.Cells.Item($row, 2).Value = $nsecode
    $ws.Cells.Item($row, 3).Value = $name
    # Force the bsecode into Excel as text (matches source "inlineStr" cells)
    # rather than letting autoconversion turn the numeric-looking string into a number.
    $ws.Cells.Item($row, 4).NumberFormat = "@"
    $ws.Cells.Item($row, 4).Value = $bsecode
    $ws.Cells.Item($row, 4).Style = "Normal"
    $ws.Cells.Item($row, 5).Value = $perChg
    $ws.Cells.Item($row, 6).Value = $close
    $ws.Cells.Item($row, 7).Value = $volume
    $ws.Cells.Item($row, 8).Value = $timeframe
    $ws.Cells.Item($row, 9).Value = $dateTime
}

Set-StockRow $ws 534 1 "ABB" "Abb India Limited" "500002" -0.86 7535.7 175164 "day" "10/09/2024 11:34:45"
Set-StockRow $ws 535 2 "HAVELLS" "Havells India Limited" "517354" 1.59 1922.45 1335776 "day" "10/09/2024 11:34:45"
Set-StockRow $ws 536 3 "PEL" "Piramal Enterprises Limited" "500302" -0.44 1058.25 486234 "day" "10/09/2024 11:34:45"
Set-StockRow $ws 537 4 "HINDALCO" "Hindalco Industries Limited" "500440" 0.14 659.5 3304741 "day" "10/09/2024 11:34:45"
Set-StockRow $ws 538 5 "IGL" "Indraprastha Gas Limited" "532514" 1.43 539.65 719416 "day" "10/09/2024 11:34:45"
Set-StockRow $ws 539 6 "ABFRL" "Aditya Birla Fashion And Retail Limited" "535755" 1.2 317.05 3612865 "day" "10/09/2024 11:34:45"
Set-StockRow $ws 540 7 "HINDCOPPER" "Hindustan Copper Limited" "513599" 1.34 313.1 3180445 "day" "10/09/2024 11:34:45"
Set-StockRow $ws 541 8 "ASHOKLEY" "Ashok Leyland Limited" "500477" 1.78 248.25 4998458 "day" "10/09/2024 11:34:45"
Set-StockRow $ws 542 9 "MOTHERSON" "Motherson Sumi Systems Limited" "517334" 1.2 189.6 10289246 "day" "10/09/2024 11:34:45"
Set-StockRow $ws 543 10 "NATIONALUM" "National Aluminium Company Limited" "532234" 2.61 175.61 8620949 "day" "10/09/2024 11:34:45"
Set-StockRow $ws 544 11 "TATASTEEL" "Tata Steel Limited" "500470" -0.03 149.42 41607844 "day" "10/09/2024 11:34:45"

